$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.406809337325484
$ws.Range("D2").Value = 0.009351139749045556
$ws.Range("E2").Value = 0.1720958756594797
$ws.Range("F2").Value = 0.9159663923201578
$ws.Range("G2").Value = 0.7963027762736203
$ws.Range("H2").Value = 0.7612044533489666
$ws.Range("L2").Value = 0.1504037954569029
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 3.128153583713527

$ws.Range("C3").Value = 0.3975481905098661
$ws.Range("D3").Value = 0.009435166500811443
$ws.Range("E3").Value = 0.167768567701188
$ws.Range("F3").Value = 0.8772088377580189
$ws.Range("G3").Value = 0.7538445135138545
$ws.Range("H3").Value = 0.7461335867586456
$ws.Range("L3").Value = 0.1462578624851503
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 3.006157524453954

$ws.Range("C4").Value = 0.3920966218257433
$ws.Range("D4").Value = 0.009488329380747729
$ws.Range("E4").Value = 0.165213771466334
$ws.Range("F4").Value = 0.8539758345155803
$ws.Range("G4").Value = 0.7282703734058202
$ws.Range("H4").Value = 0.7373056249113006
$ws.Range("L4").Value = 0.1438038777470965
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 2.933175414339928

$ws.Range("C5").Value = 0.3899339996493438
$ws.Range("D5").Value = 0.009510388064757791
$ws.Range("E5").Value = 0.1641983007239354
$ws.Range("F5").Value = 0.8446496412523885
$ws.Range("G5").Value = 0.7179726945568348
$ws.Range("H5").Value = 0.7338150516635267
$ws.Range("L5").Value = 0.1428268243801014
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 2.90391688686924

$ws.Range("C6").Value = 0.3895784543456386
$ws.Range("D6").Value = 0.009514074699971564
$ws.Range("E6").Value = 0.1640312286104617
$ws.Range("F6").Value = 0.8431095637789809
$ws.Range("G6").Value = 0.7162702446706817
$ws.Range("H6").Value = 0.7332418969717196
$ws.Range("L6").Value = 0.1426659702312065
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 2.899087614875384

$ws.Range("C7").Value = 0.3920672174581625
$ws.Range("D7").Value = 0.009488625275758444
$ws.Range("E7").Value = 0.1651999727869686
$ws.Range("F7").Value = 0.8538494861101213
$ws.Range("G7").Value = 0.7281309940445908
$ws.Range("H7").Value = 0.7372581171935622
$ws.Range("L7").Value = 0.1437906079780333
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 2.932778872918107

$ws.Range("C8").Value = 0.4035672883946972
$ws.Range("D8").Value = 0.009379786018809844
$ws.Range("E8").Value = 0.1705825686378546
$ws.Range("F8").Value = 0.9024853816899849
$ws.Range("G8").Value = 0.7815599795247579
$ws.Range("H8").Value = 0.7559195968969448
$ws.Range("L8").Value = 0.1489552107992722
$ws.Range("N8").Value = 2.766433886209654
$ws.Range("O8").Value = 3.085689059124093

$ws.Range("C9").Value = 0.4279886984083419
$ws.Range("D9").Value = 0.009178828835495256
$ws.Range("E9").Value = 0.1819524292491792
$ws.Range("F9").Value = 1.002364138748504
$ws.Range("G9").Value = 0.8902974908078249
$ws.Range("H9").Value = 0.7959006314800376
$ws.Range("L9").Value = 0.1598142602100125
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 3.400899602077345

$ws.Range("C10").Value = 0.4470827431900943
$ws.Range("D10").Value = 0.009038814223037761
$ws.Range("E10").Value = 0.1908088464604418
$ws.Range("F10").Value = 1.078536330038474
$ws.Range("G10").Value = 0.9726574500498657
$ws.Range("H10").Value = 0.8273555155266763
$ws.Range("L10").Value = 0.16824506866962
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 3.641995857184952

$ws.Range("C11").Value = 0.4560218970184167
$ws.Range("D11").Value = 0.008976781068605
$ws.Range("E11").Value = 0.1949485316877713
$ws.Range("F11").Value = 1.113805546471013
$ws.Range("G11").Value = 1.010673884672968
$ws.Range("H11").Value = 0.8421208982347537
$ws.Range("L11").Value = 0.1721803059860036
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 3.75377647496191

$ws.Range("C12").Value = 0.4594434993398977
$ws.Range("D12").Value = 0.008953530418285949
$ws.Range("E12").Value = 0.1965321595870648
$ws.Range("F12").Value = 1.127250626168532
$ws.Range("G12").Value = 1.025149684591923
$ws.Range("H12").Value = 0.8477780284605387
$ws.Range("L12").Value = 0.1736849660906188
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 3.796409758596212

$ws.Range("C13").Value = 0.4587049688618947
$ws.Range("D13").Value = 0.008958527169437236
$ws.Range("E13").Value = 0.196190383286627
$ws.Range("F13").Value = 1.12435100373186
$ws.Range("G13").Value = 1.022028502582799
$ws.Range("H13").Value = 0.8465567350012009
$ws.Range("L13").Value = 0.1733602655908442
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 3.787214352518163

$ws.Range("C14").Value = 0.4563026612448198
$ws.Range("D14").Value = 0.008974863407454681
$ws.Range("E14").Value = 0.1950784962148759
$ws.Range("F14").Value = 1.11490988676033
$ws.Range("G14").Value = 1.011863214402609
$ws.Range("H14").Value = 0.8425849940171588
$ws.Range("L14").Value = 0.1723038047457663
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 3.757277830442263

$ws.Range("C15").Value = 0.4548359408504723
$ws.Range("D15").Value = 0.008984901112563559
$ws.Range("E15").Value = 0.194399521746945
$ws.Range("F15").Value = 1.109138586840288
$ws.Range("G15").Value = 1.00564709354552
$ws.Range("H15").Value = 0.8401607622488143
$ws.Range("L15").Value = 0.1716585791815533
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 3.738980538580222

$ws.Range("C16").Value = 0.4465036542897565
$ws.Range("D16").Value = 0.009042901745663023
$ws.Range("E16").Value = 0.1905405456509044
$ws.Range("F16").Value = 1.0762438851241
$ws.Range("G16").Value = 0.9701841184885609
$ws.Range("H16").Value = 0.8263997614479308
$ws.Range("L16").Value = 0.1679899118500998
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 3.634733230999814

$ws.Range("C17").Value = 0.4414570045306334
$ws.Range("D17").Value = 0.00907890930486599
$ws.Range("E17").Value = 0.1882016369688557
$ws.Range("F17").Value = 1.056222723677024
$ws.Range("G17").Value = 0.9485701779533144
$ws.Range("H17").Value = 0.818074847965363
$ws.Range("L17").Value = 0.1657649834153574
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 3.571321180802727

$ws.Range("C18").Value = 0.4385781272713132
$ws.Range("D18").Value = 0.009099776085976785
$ws.Range("E18").Value = 0.1868667881919635
$ws.Range("F18").Value = 1.044765229239161
$ws.Range("G18").Value = 0.9361901856157431
$ws.Range("H18").Value = 0.813329525841624
$ws.Range("L18").Value = 0.1644946724680381
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 3.535046248280196

$ws.Range("C19").Value = 0.4376074763868019
$ws.Range("D19").Value = 0.009106868008980529
$ws.Range("E19").Value = 0.186416620273782
$ws.Range("F19").Value = 1.040895890022568
$ws.Range("G19").Value = 0.9320074057541774
$ws.Range("H19").Value = 0.811730212508138
$ws.Range("L19").Value = 0.164066180343994
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 3.522798133886681

$ws.Range("C20").Value = 0.4419917625539256
$ws.Range("D20").Value = 0.009075060065984131
$ws.Range("E20").Value = 0.188449537801354
$ws.Range("F20").Value = 1.058347988721067
$ws.Range("G20").Value = 0.9508656557133861
$ws.Range("H20").Value = 0.8189566030620767
$ws.Range("L20").Value = 0.1660008564073507
$ws.Range("N20").Value = 3.957806003280837
$ws.Range("O20").Value = 3.578050996947979

$ws.Range("C21").Value = 0.4570072845097002
$ws.Range("D21").Value = 0.008970058535879621
$ws.Range("E21").Value = 0.195404649025015
$ws.Range("F21").Value = 1.117680541478578
$ws.Range("G21").Value = 1.014846835011781
$ws.Range("H21").Value = 0.8437498034405735
$ws.Range("L21").Value = 0.1726137194532527
$ws.Range("N21").Value = 4.391158149571254
$ws.Range("O21").Value = 3.766062637682523

$ws.Range("C22").Value = 0.4670338348652194
$ws.Range("D22").Value = 0.008902833277983024
$ws.Range("E22").Value = 0.2000436238092007
$ws.Range("F22").Value = 1.156979085665597
$ws.Range("G22").Value = 1.05712767971977
$ws.Range("H22").Value = 0.8603372093284065
$ws.Range("L22").Value = 0.1770199983058518
$ws.Range("N22").Value = 4.673791817957863
$ws.Range("O22").Value = 3.890714148656627

$ws.Range("C23").Value = 0.4616629389616662
$ws.Range("D23").Value = 0.008938584191028331
$ws.Range("E23").Value = 0.1975591428688332
$ws.Range("F23").Value = 1.135956844762447
$ws.Range("G23").Value = 1.034518794695003
$ws.Range("H23").Value = 0.851449040292664
$ws.Range("L23").Value = 0.1746605331654081
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 3.824022297306726

$ws.Range("C24").Value = 0.4417499283512996
$ws.Range("D24").Value = 0.009076799791260637
$ws.Range("E24").Value = 0.1883374312510355
$ws.Range("F24").Value = 1.057386991536674
$ws.Range("G24").Value = 0.9498277265160766
$ws.Range("H24").Value = 0.818557834561858
$ws.Range("L24").Value = 0.1658941907312084
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 3.575007881843646

$ws.Range("C25").Value = 0.4211806766275288
$ws.Range("D25").Value = 0.009231855463464278
$ws.Range("E25").Value = 0.1787887158458972
$ws.Range("F25").Value = 0.9748574239048935
$ws.Range("G25").Value = 0.8604512084783096
$ws.Range("H25").Value = 0.7847205061395357
$ws.Range("L25").Value = 0.1567976055938942
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 3.313967518965057
